# Add a new "UK" test-data sheet, modelled on the existing "Netherlands"
# sheet, and populate it with the UK market values.

$wb = $excel.ActiveWorkbook

# Duplicate the Netherlands sheet (keeps column widths, styles, merged
# cells, etc. identical) and place the copy right after it, i.e. as the
# new last tab.
$netherlands = $wb.Worksheets.Item("Netherlands")
$netherlands.Copy($null, $netherlands)

$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# The copied row 2 inherited Netherlands' taller (wrapped) row height;
# restore the default auto height used by the other country sheets.
$uk.Rows.Item(2).AutoFit()

# Fill in the UK-specific values (code first, then market name, to match
# shared-string insertion order).
$uk.Range("B4").Value = "NGC-2741/T3396"
$uk.Range("B2").Value = "UK Market"

# Leave the new sheet active with B4 selected, matching the saved view.
$uk.Range("B4").Select() | Out-Null
